$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1100.8
$ws.Range("I6").Value = 114.17647
$ws.Range("J6").Value = 6691.6665
$ws.Range("K6").Value = 342.52941
$ws.Range("L6").Value = 20074.9995
$ws.Range("M6").Value = -230.52941
$ws.Range("N6").Value = -20298.9995
$ws.Range("H9").Value = 63.5
$ws.Range("I9").Value = 56.2
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 56.2
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = 112.8
$ws.Range("N9").Value = -438
$ws.Range("H19").Value = 2098.2222
$ws.Range("I19").Value = 2454.4285
$ws.Range("J19").Value = 1871.5454
$ws.Range("K19").Value = 2454.4285
$ws.Range("L19").Value = 1871.5454
$ws.Range("M19").Value = -2279.4285
$ws.Range("N19").Value = -2221.5454
$ws.Range("H33").Value = 958334.5
$ws.Range("I33").Value = 1077838.8
$ws.Range("J33").Value = 2300
$ws.Range("K33").Value = 1077838.8
$ws.Range("L33").Value = 2300
$ws.Range("M33").Value = -1077609.8
$ws.Range("N33").Value = -2758
$ws.Range("H39").Value = 433
$ws.Range("I39").Value = 225.53847
$ws.Range("J39").Value = 1332
$ws.Range("K39").Value = 676.61541
$ws.Range("L39").Value = 3996
$ws.Range("M39").Value = -380.61541
$ws.Range("N39").Value = -4588
$ws.Range("H75").Value = 25314
$ws.Range("J75").Value = 25314
$ws.Range("L75").Value = 25314
$ws.Range("N75").Value = -27186
$ws.Range("H78").Value = 25314
$ws.Range("J78").Value = 25314
$ws.Range("L78").Value = 75942
$ws.Range("N78").Value = -85302
$ws.Range("H106").Value = 10021.333
$ws.Range("I106").Value = 2935
$ws.Range("J106").Value = 24194
$ws.Range("K106").Value = 2935
$ws.Range("L106").Value = 24194
$ws.Range("M106").Value = -2304
$ws.Range("N106").Value = -25456
$ws.Range("H132").Value = 1310223
$ws.Range("I132").Value = 1434936.4
$ws.Range("K132").Value = 4304809.199999999
$ws.Range("M132").Value = -4302279.199999999
$ws.Range("H137").Value = 69722
$ws.Range("I137").Value = 812
$ws.Range("J137").Value = 86949.5
$ws.Range("K137").Value = 2436
$ws.Range("L137").Value = 260848.5
$ws.Range("M137").Value = 114
$ws.Range("N137").Value = -265948.5
$ws.Range("H138").Value = 5496.9287
$ws.Range("I138").Value = 3076.5
$ws.Range("J138").Value = 7312.25
$ws.Range("K138").Value = 9229.5
$ws.Range("L138").Value = 21936.75
$ws.Range("M138").Value = -4089.5
$ws.Range("N138").Value = -32216.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 124807.42
$ws.Range("I32").Value = 132305.34
$ws.Range("J32").Value = 93452.45
$ws.Range("K32").Value = 132305.34
$ws.Range("L32").Value = 93452.45
$ws.Range("M32").Value = -132018.34
$ws.Range("N32").Value = -94026.45
$ws.Range("H97").Value = 1327
$ws.Range("I97").Value = 1014.5217
$ws.Range("J97").Value = 3123.75
$ws.Range("K97").Value = 1014.5217
$ws.Range("L97").Value = 3123.75
$ws.Range("M97").Value = -518.5217
$ws.Range("N97").Value = -4115.75
$ws.Range("H132").Value = 12736.429
$ws.Range("I132").Value = 13769.521
$ws.Range("K132").Value = 41308.563
$ws.Range("M132").Value = -38778.563

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 20068.223
$ws.Range("I82").Value = 10204.667
$ws.Range("J82").Value = 25000
$ws.Range("K82").Value = 10204.667
$ws.Range("L82").Value = 25000
$ws.Range("M82").Value = -9821.666999999999
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 20068.223
$ws.Range("I85").Value = 10204.667
$ws.Range("J85").Value = 25000
$ws.Range("K85").Value = 10204.667
$ws.Range("L85").Value = 25000
$ws.Range("M85").Value = -8878.666999999999
$ws.Range("N85").Value = -27652
$ws.Range("H107").Value = 7860.067
$ws.Range("I107").Value = 7561.6924
$ws.Range("J107").Value = 9799.5
$ws.Range("K107").Value = 7561.6924
$ws.Range("L107").Value = 9799.5
$ws.Range("M107").Value = -5641.6924
$ws.Range("N107").Value = -13639.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3185.8794
$ws.Range("I31").Value = 2507.2307
$ws.Range("J31").Value = 3737.2812
$ws.Range("K31").Value = 2507.2307
$ws.Range("L31").Value = 3737.2812
$ws.Range("M31").Value = -2212.2307
$ws.Range("N31").Value = -4327.281199999999
$ws.Range("H34").Value = 3185.8794
$ws.Range("I34").Value = 2507.2307
$ws.Range("J34").Value = 3737.2812
$ws.Range("K34").Value = 2507.2307
$ws.Range("L34").Value = 3737.2812
$ws.Range("M34").Value = -2305.2307
$ws.Range("N34").Value = -4141.281199999999
$ws.Range("H93").Value = 16393.125
$ws.Range("I93").Value = 7999.8
$ws.Range("K93").Value = 7999.8
$ws.Range("M93").Value = -6127.8
$ws.Range("H132").Value = 3119.8125
$ws.Range("I132").Value = 2961.8333
$ws.Range("J132").Value = 5489.5
$ws.Range("K132").Value = 8885.499899999999
$ws.Range("L132").Value = 16468.5
$ws.Range("M132").Value = -6355.499899999999
$ws.Range("N132").Value = -21528.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4170754.5
$ws.Range("J68").Value = 12504139
$ws.Range("L68").Value = 37512417
$ws.Range("N68").Value = -37514039
$ws.Range("H71").Value = 4170754.5
$ws.Range("J71").Value = 12504139
$ws.Range("L71").Value = 112537251
$ws.Range("N71").Value = -112545363
$ws.Range("H86").Value = 675
$ws.Range("I86").Value = 400
$ws.Range("J86").Value = 950
$ws.Range("K86").Value = 1200
$ws.Range("L86").Value = 2850
$ws.Range("M86").Value = -14
$ws.Range("N86").Value = -5222
$ws.Range("H89").Value = 675
$ws.Range("I89").Value = 400
$ws.Range("J89").Value = 950
$ws.Range("K89").Value = 3600
$ws.Range("L89").Value = 8550
$ws.Range("M89").Value = 2328
$ws.Range("N89").Value = -20406
$ws.Range("H131").Value = 111581.4
$ws.Range("J131").Value = 232713.3
$ws.Range("L131").Value = 698139.8999999999
$ws.Range("N131").Value = -708219.8999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 32562.75
$ws.Range("J93").Value = 32562.75
$ws.Range("L93").Value = 32562.75
$ws.Range("N93").Value = -36306.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 24795.111
$ws.Range("I46").Value = 53995.625
$ws.Range("J46").Value = 1434.7
$ws.Range("K46").Value = 53995.625
$ws.Range("L46").Value = 1434.7
$ws.Range("M46").Value = -53807.625
$ws.Range("N46").Value = -1810.7
$ws.Range("H61").Value = 30335616
$ws.Range("I61").Value = 33335844
$ws.Range("K61").Value = 33335844
$ws.Range("M61").Value = -33335642
$ws.Range("H82").Value = 1897.6
$ws.Range("I82").Value = 999
$ws.Range("J82").Value = 2122.25
$ws.Range("K82").Value = 999
$ws.Range("L82").Value = 2122.25
$ws.Range("M82").Value = -638
$ws.Range("N82").Value = -2844.25
$ws.Range("H85").Value = 1897.6
$ws.Range("I85").Value = 999
$ws.Range("J85").Value = 2122.25
$ws.Range("K85").Value = 999
$ws.Range("L85").Value = 2122.25
$ws.Range("M85").Value = 249
$ws.Range("N85").Value = -4618.25
$ws.Range("H93").Value = 32007.092
$ws.Range("I93").Value = 1625.1666
$ws.Range("J93").Value = 68465.39999999999
$ws.Range("K93").Value = 1625.1666
$ws.Range("L93").Value = 68465.39999999999
$ws.Range("M93").Value = -377.1666
$ws.Range("N93").Value = -70961.39999999999
$ws.Range("H113").Value = 30335616
$ws.Range("I113").Value = 33335844
$ws.Range("K113").Value = 33335844
$ws.Range("M113").Value = -33333674

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4809.294
$ws.Range("I136").Value = 5339.857
$ws.Range("K136").Value = 16019.571
$ws.Range("M136").Value = -13469.571
